$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 263 (shifts existing rows 263:351 down to 264:352,
# bringing row formatting/styles along with it, matching dimension A1:T352).
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row with the new weekly price-point data.
$ws.Cells.Item(263, 1).Value = 7
$ws.Cells.Item(263, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(263, 3).Value = "Ñuble"
$ws.Cells.Item(263, 4).Value = 45135
$ws.Cells.Item(263, 5).Value = 16
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100101
$ws.Cells.Item(263, 8).Value = "Berries"
$ws.Cells.Item(263, 9).Value = 100101007
$ws.Cells.Item(263, 10).Value = "Kiwi"
$ws.Cells.Item(263, 11).Value = "Hayward"
$ws.Cells.Item(263, 12).Value = "Primera"
$ws.Cells.Item(263, 13).Value = 60
$ws.Cells.Item(263, 14).Value = 14000
$ws.Cells.Item(263, 15).Value = 14000
$ws.Cells.Item(263, 16).Value = 14000
$ws.Cells.Item(263, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(263, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(263, 19).Value = 778
$ws.Cells.Item(263, 20).Value = 18
